# "Agregado el flujo 3"
# - Rename "Hoja 3" -> "validarCat" and repurpose it as a product list.
# - Add a brand-new sheet "validarHome" with another product list (white-filled rows).
# - Clear the old xpath/nombre helper columns on "productos".

$wb = $excel.ActiveWorkbook

# --- productos: drop the xpath/nombre helper columns (B:C), keep column A ---
$wsProductos = $wb.Worksheets.Item("productos")
$wsProductos.Range("B1:C4").Clear()

# --- Hoja 3 -> validarCat: replace the step-by-step flow table with a product list ---
$wsCat = $wb.Worksheets.Item("Hoja 3")
$wsCat.Name = "validarCat"
$wsCat.Range("A1:C6").Clear()

$wsCat.Columns.Item(1).ColumnWidth = 25.75

$wsCat.Range("A1").Value = "Producto"
$wsCat.Range("A2").Value = "Encendedor Zippo Calavera Naipe"
$wsCat.Range("A3").Value = "Estuche Zippo Cuero Cafe / Correa " + [char]0x2013 + " Cod Lplb"
$wsCat.Range("A4").Value = "Inserto Zippo Llama Amarilla"

$wsProductos.Range("A1").Copy()
$wsCat.Range("A1:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- new sheet: validarHome ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsHome = $wb.Worksheets.Add($null, $lastSheet)
$wsHome.Name = "validarHome"

$wsHome.Range("A1").Value = "Texto"
$wsHome.Range("A2").Value = "Combustible Para Encendedor Zippo 4oz " + [char]0x2013 + " Cod 3141laex"
$wsHome.Range("A3").Value = "Dispensador Guarda Combustible Zippo " + [char]0x2013 + " Cod 121503"
$wsHome.Range("A4").Value = "Inserto Encendedor Zippo Butano Llama Sencilla Cod 65826"
$wsHome.Range("A5").Value = "Inserto Encendedor Zippo Butano Llama Doble Cod 65827"
$wsHome.Range("A6").Value = "Mecha Para Encendedor Zippo " + [char]0x2013 + " Cod 2425"

$wsProductos.Range("A1").Copy()
$wsHome.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsHome.Range("A2:A6").Interior.Color = 16777215
